$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A (shifts every existing column right by one)
$ws.Columns.Item(1).Insert()

# New column header ("Match ID") in the visible header row
$ws.Range("A2").Value() = "Match ID"
$ws.Range("A2").Font.Bold = $true

# Blank styled cell in the hidden spacer row
$ws.Range("A3").Font.Bold = $true

# Data rows: constant Match ID value of 33, bold style (matches header font)
for ($r = 4; $r -le 19; $r++) {
    $ws.Cells.Item($r, 1).Value() = 33
    $ws.Cells.Item($r, 1).Font.Bold = $true
}

# Hidden totals row keeps the default (non-bold) style
$ws.Cells.Item(20, 1).Value() = 33
$ws.Rows.Item(20).AutoFit()

# Update selection to match the new layout
$null = $ws.Range("A2:A19").Select()

